$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new blank columns before the existing "fantasy points"
# column (E), shifting rec_yds/rec_td/fumbles headers data untouched and
# pushing "fantasy points" from column E to column G.
$ws.Range("E1").EntireColumn.Insert()
$ws.Range("F1").EntireColumn.Insert()

# Copy the header style (bold, centered, bordered) from the neighboring
# "fumbles" header cell into the two new header cells.
$ws.Range("D1").Copy($ws.Range("E1"))
$ws.Range("D1").Copy($ws.Range("F1"))

# Label the new header cells.
$ws.Range("E1").Value = "height"
$ws.Range("F1").Value = "weight"

# Fill the new columns with a constant value for every player/data row
# (row 1 is the header, data starts on row 2).
$lastRow = $ws.UsedRange.Rows.Count
for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 5).Value = 6.416666666666667
    $ws.Cells.Item($r, 6).Value = 253
}
